$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 388
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 46062
}
